$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.668.37"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "3.802.38"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'613.95"
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("D6").Value = "'177.21"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("D7").Value = "3.801.83"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").Value = "'6.50"
$ws.Range("E11").Value = "  +2.70%  "
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "'39.79"
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "4.436.88"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").Value = "3.806.54"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "69.738.51"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("D20").Value = "'16.67"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "'508.39"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").Value = "'9.62"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").Value = "'86.36"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("E26").Value = "  +3.69%  "
$ws.Range("D27").Value = "'12.72"
$ws.Range("E27").Value = "  -3.58%  "
$ws.Range("E28").Value = "  -5.59%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").Value = "'8.07"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").Value = "'31.49"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("D37").Value = "'6.12"
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("D38").Value = "'0.141"
$ws.Range("E38").Value = "  +6.52%  "
$ws.Range("D39").Value = "'479.88"
$ws.Range("E39").Value = "  +12.96%  "
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("E41").Value = "  +5.35%  "
$ws.Range("D42").Value = "'2.07"
$ws.Range("E42").Value = "  -2.53%  "
$ws.Range("D43").Value = "'49.75"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").Value = "'44.17"
$ws.Range("E44").Value = "  -2.64%  "
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("D46").Value = "2.933.84"
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "'27.27"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "'139.26"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  -3.78%  "